$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.138.41'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.926.85'
$ws.Range('E3').Value = '  +3.29%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '470.69'
$ws.Range('E5').Value = '  +8.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.01'
$ws.Range('E6').Value = '  +3.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.620'
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.728'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  +6.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000338'
$ws.Range('E11').Value = '  +4.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.15'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.539.95'
$ws.Range('E13').Value = '  +3.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.95'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.924.32'
$ws.Range('E16').Value = '  +3.67%  '
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.79'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.380.88'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '429.81'
$ws.Range('E21').Value = '  +4.15%  '
$ws.Range('E22').Value = '  +2.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.51'
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.41'
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.55'
$ws.Range('E25').Value = '  +5.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '38.38'
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.74'
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.18'
$ws.Range('E28').Value = '  +3.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.57'
$ws.Range('E29').Value = '  -2.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '727.52'
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.132'
$ws.Range('E31').Value = '  -3.40%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.52'
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.80'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.63'
$ws.Range('E34').Value = '  +1.90%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.88'
$ws.Range('E35').Value = '  +3.69%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.151'
$ws.Range('E36').Value = '  +0.45%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0784'
$ws.Range('E38').Value = '  +10.82%  '
$ws.Range('E39').Value = '  -5.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0476'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('E41').Value = '  +3.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.58'
$ws.Range('E42').Value = '  -5.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.140'
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.336'
$ws.Range('E45').Value = '  +3.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('E46').Value = '  +4.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.18'
$ws.Range('E47').Value = '  +4.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.42'
$ws.Range('E48').Value = '  +2.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '146.32'
$ws.Range('E49').Value = '  +3.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.16'
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('E51').Value = '  +1.42%  '
